$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 276, pushing existing rows 276.. down to 278..
$ws.Rows("276:277").Insert()

# Populate new row 276
$ws.Range("A276").Value = 8
$ws.Range("B276").Value = "Terminal La Palmera de La Serena"
$ws.Range("C276").Value = "Coquimbo"
$ws.Range("D276").Value = 44946
$ws.Range("E276").Value = 4
$ws.Range("F276").Value = "Fruta"
$ws.Range("G276").Value = 100103
$ws.Range("H276").Value = "Frutos de hueso (carozo)"
$ws.Range("I276").Value = 100103002
$ws.Range("J276").Value = "Ciruela"
$ws.Range("K276").Value = "Black Amber"
$ws.Range("L276").Value = "Primera"
$ws.Range("M276").Value = 10
$ws.Range("N276").Value = 330000
$ws.Range("O276").Value = 340000
$ws.Range("P276").Value = 335000
$ws.Range("Q276").Value = "`$/bins (450 kilos)"
$ws.Range("R276").Value = "Provincia de Curicó"
$ws.Range("S276").Value = 744
$ws.Range("T276").Value = 450

# Populate new row 277
$ws.Range("A277").Value = 8
$ws.Range("B277").Value = "Terminal La Palmera de La Serena"
$ws.Range("C277").Value = "Coquimbo"
$ws.Range("D277").Value = 44946
$ws.Range("E277").Value = 4
$ws.Range("F277").Value = "Fruta"
$ws.Range("G277").Value = 100103
$ws.Range("H277").Value = "Frutos de hueso (carozo)"
$ws.Range("I277").Value = 100103002
$ws.Range("J277").Value = "Ciruela"
$ws.Range("K277").Value = "Black Amber"
$ws.Range("L277").Value = "Segunda"
$ws.Range("M277").Value = 10
$ws.Range("N277").Value = 290000
$ws.Range("O277").Value = 300000
$ws.Range("P277").Value = 295000
$ws.Range("Q277").Value = "`$/bins (450 kilos)"
$ws.Range("R277").Value = "Provincia de Curicó"
$ws.Range("S277").Value = 656
$ws.Range("T277").Value = 450
